$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 (Experimental row) needs to contain the literal TEXT "false" (not a
# boolean). Assigning the string directly gets auto-coerced to a Boolean by
# Excel's smart literal parsing, so build it as a text formula first, then
# convert it to a plain value via copy / paste-special-values (this keeps
# the cell's existing style intact instead of picking up a new number
# format / quote-prefix).
$ws.Range("B7").Formula = '="false"'
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)

$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"
$ws.Range("B17").Value = "Validation status of health measurements"
